$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two shared strings used in G1/H1 header cells
$ws.Range("G1").Value = "DP_口頭報告"
$ws.Range("H1").Value = "DV_口頭報告"

# Column width changes (stored OOXML "width" = ColumnWidth + 5/6 in this
# engine's pixel-grid model, so back the COM value off by 5/6 to land the
# saved width attribute on the intended target)
$ws.Columns.Item(3).ColumnWidth = 10.41666666666667
$ws.Columns.Item(7).ColumnWidth = 12.04666666666667
$ws.Columns.Item(8).ColumnWidth = 11.66666666666667
